# 5.5.13 Real Property-Monthly Reviews-org_results.xlsx
# Update "Resource Utilization" sheet (B2/B3 -> 0) and truncate the
# "Activity Times" sheet down to just the header + the summary "Process"
# row, zeroing out the numeric columns on that row and deleting the
# now-obsolete "Activity Step"/"Stop" detail rows (3-6).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: Resource Utilization ---
$wsResource = $wb.Worksheets.Item("Resource Utilization")
$wsResource.Range("B2").Value = 0
$wsResource.Range("B3").Value = 0

# --- Sheet 2: Activity Times ---
$wsActivity = $wb.Worksheets.Item("Activity Times")

# Zero out the aggregate numbers on row 2 (C2:G2); H2:K2 are already 0.
$wsActivity.Range("C2:G2").Value = 0

# Remove rows 3-6 (the per-step detail rows) entirely.
$wsActivity.Rows("3:6").Delete()
